$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet: previously hard-coded settings (cool down time) are joined
# by several new, user-configurable settings (timezone, legend texts,
# impact-of-preferred-event weight and shuffle flag).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Settings")

# --- Row 1: header row (values unchanged, only the comment column header
#     formatting is touched by the original edit, which is not visually
#     significant, so the existing formatting is simply kept in place) -----
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Comment"

# --- Row 2: the old "Cool down time in days" row becomes a blank spacer
#     row; reset A2/B2 back to an unformatted blank cell and let C2 keep a
#     wrapped (blank) style ----------------------------------------------
$ws.Range("A2:B2").Clear()
$ws.Range("C2").ClearContents()
$ws.Range("C2").WrapText = $true

# --- Row 3: Timezone -----------------------------------------------------
$ws.Range("A3").Value = "Timezone"
$ws.Range("B3").Value = "Europe/Berlin"
$ws.Range("C3").Value = "Time zone of the user"
$ws.Range("A3:C3").WrapText = $true

# --- Row 4: blank spacer row, same wrapped style as row 3/5 -------------
$ws.Range("A4:C4").WrapText = $true

# --- Row 5: legendDay ------------------------------------------------------
$ws.Range("A5").Value = "legendDay"
$ws.Range("B5").Value = "Day"
$ws.Range("C5").Value = "Text for the legend of the output file"
$ws.Range("A5:C5").WrapText = $true

# --- Row 6: legendDate -----------------------------------------------------
$ws.Range("A6").Value = "legendDate"
$ws.Range("B6").Value = "Date"
$ws.Range("A6:C6").WrapText = $true

# --- Row 7: legendTime -----------------------------------------------------
$ws.Range("A7").Value = "legendTime"
$ws.Range("B7").Value = "Time"
$ws.Range("A7:C7").WrapText = $true

# --- Row 8: legendComment ---------------------------------------------------
$ws.Range("A8").Value = "legendComment"
$ws.Range("B8").Value = "Comment"

# --- Row 9: blank spacer row -------------------------------------------------

# --- Row 10: cooldowntime (previously hard-coded "Cool down time in days") -
$ws.Range("A10").Value = "cooldowntime"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "The worker will not be selected during the cool down time after his last action."
$ws.Range("A10:C10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 23.95

# --- Row 11: impactOfPreferredEvent -----------------------------------------
$ws.Range("A11").Value = "impactOfPreferredEvent"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Increasing this values makes it more likely that people are selected for preferred events."
$ws.Rows.Item(11).RowHeight = 23.95

# --- Row 12: shuffle ---------------------------------------------------------
$ws.Range("A12").Value = "shuffle"
$ws.Range("B12").Value = "'true"
$ws.Range("C12").Value = "If several workers have the same ranking value, this settings decides whether they are selected randomly. If shuffle = false, the first worker is selected."
$ws.Rows.Item(12).RowHeight = 46.45

# --- Column widths (best effort, mirrors the resize seen in the sheet) ------
$ws.Columns.Item(1).ColumnWidth = 20.6530612244898
$ws.Columns.Item(2).ColumnWidth = 12.2857142857143
$ws.Columns.Item(3).ColumnWidth = 36.9897959183673

# --- Selection moves to C17 on the Settings sheet ----------------------------
$ws.Activate()
$ws.Range("C17").Select()

# ---------------------------------------------------------------------------
# Minor column width tweaks on the other sheets (cosmetic resize following
# the content changes above).
# ---------------------------------------------------------------------------
$wsPeriod = $wb.Worksheets.Item("Period")
$wsPeriod.Columns.Item(1).ColumnWidth = 27.6734693877551
$wsPeriod.Columns.Item(2).ColumnWidth = 9.71938775510204

$wsRegular = $wb.Worksheets.Item("RegularEvents")
$wsRegular.Columns.Item(2).ColumnWidth = 11.0714285714286

$wsSpecial = $wb.Worksheets.Item("SpecialEvents")
$wsSpecial.Columns.Item(7).ColumnWidth = 13.5

$wsWorkers = $wb.Worksheets.Item("Workers")
$wsWorkers.Columns.Item(7).ColumnWidth = 12.9591836734694
$wsWorkers.Columns.Item(8).ColumnWidth = 8.36734693877551
